# Applies the "Updated cryptos list" refresh: new Price (D) / Volume(1h) (E)
# figures scraped for each coin, plus two rank positions (32/33 and 43/44)
# where the scraped coin order swapped between runs so Coin/Link/Price/Volume
# moved to the other row of the pair while the rank index (column A) stayed put.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.222.12'
$ws.Range("E2").Value = '  +0.51%  '
$ws.Range("D3").Value = '2.475.45'
$ws.Range("E3").Value = '  +0.55%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").Value = '''518.09'
$ws.Range("E5").Value = '  +0.29%  '
$ws.Range("D6").Value = '''131.06'
$ws.Range("E6").Value = '  +0.28%  '
$ws.Range("D7").Value = '''0.997'
$ws.Range("E7").Value = '  -0.36%  '
$ws.Range("D8").Value = '''0.553'
$ws.Range("E8").Value = '  -0.69%  '
$ws.Range("D9").Value = '2.506.07'
$ws.Range("E9").Value = '  +1.70%  '
$ws.Range("E10").Value = '  -1.99%  '
$ws.Range("E11").Value = '  -0.10%  '
$ws.Range("E12").Value = '  -2.56%  '
$ws.Range("D13").Value = '''0.329'
$ws.Range("E13").Value = '  -2.98%  '
$ws.Range("D14").Value = '2.919.67'
$ws.Range("E14").Value = '  +0.57%  '
$ws.Range("D15").Value = '58.087.02'
$ws.Range("E15").Value = '  +0.31%  '
$ws.Range("D16").Value = '''21.97'
$ws.Range("E16").Value = '  -1.22%  '
$ws.Range("E17").Value = '  -1.22%  '
$ws.Range("D18").Value = '2.496.27'
$ws.Range("E18").Value = '  +1.06%  '
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("D20").Value = '''319.87'
$ws.Range("E20").Value = '  +0.55%  '
$ws.Range("E21").Value = '  -0.27%  '
$ws.Range("D22").Value = '''0.995'
$ws.Range("E22").Value = '  -0.43%  '
$ws.Range("E23").Value = '  +4.48%  '
$ws.Range("D24").Value = '''64.21'
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("E25").Value = '  -2.01%  '
$ws.Range("D26").Value = '''0.993'
$ws.Range("E26").Value = '  -0.69%  '
$ws.Range("D27").Value = '''0.159'
$ws.Range("E27").Value = '  +0.49%  '
$ws.Range("D28").Value = '''7.30'
$ws.Range("E28").Value = '  +0.15%  '
$ws.Range("E29").Value = '  +0.75%  '
$ws.Range("D30").Value = '''167.55'
$ws.Range("E30").Value = '  +1.27%  '
$ws.Range("E31").Value = '  +0.86%  '
# Row 32: coin identity + price/volume swapped in with its pair row
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").Value = '''1.17'
$ws.Range("E32").Value = '  +1.34%  '
# Row 33: coin identity + price/volume swapped in with its pair row
$ws.Range("B33").Value = 'Aptos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D33").Value = '''6.21'
$ws.Range("E33").Value = '  -0.30%  '
$ws.Range("D35").Value = '''0.995'
$ws.Range("E35").Value = '  -0.59%  '
$ws.Range("D36").Value = '''17.97'
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("E37").Value = '  -2.62%  '
$ws.Range("E38").Value = '  -1.08%  '
$ws.Range("D39").Value = '''36.70'
$ws.Range("E39").Value = '  +0.79%  '
$ws.Range("E40").Value = '  -1.09%  '
$ws.Range("D41").Value = '''0.764'
$ws.Range("E41").Value = '  -2.56%  '
$ws.Range("D42").Value = '''274.18'
$ws.Range("E42").Value = '  +1.58%  '
# Row 43: coin identity + price/volume swapped in with its pair row
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").Value = '''3.41'
$ws.Range("E43").Value = '  -0.78%  '
# Row 44: coin identity + price/volume swapped in with its pair row
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '''4.95'
$ws.Range("E44").Value = '  +0.60%  '
$ws.Range("E45").Value = '  +0.07%  '
$ws.Range("D46").Value = '''0.0917'
$ws.Range("E46").Value = '  +1.41%  '
$ws.Range("D47").Value = '''120.36'
$ws.Range("E47").Value = '  -4.84%  '
$ws.Range("E48").Value = '  +2.38%  '
$ws.Range("D49").Value = '''17.62'
$ws.Range("E49").Value = '  -0.08%  '
$ws.Range("E50").Value = '  +1.46%  '
$ws.Range("D51").Value = '''16.68'
$ws.Range("E51").Value = '  -0.56%  '
